$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("H4").Value = 5.5
$ws.Range("I4").Value = 4.5
$ws.Range("Q4").Value = 113
$ws.Range("R4").Value = 104
$ws.Range("T4").Value = 0.1131736743200052
$ws.Range("U4").Value = 0.1230816648328576
$ws.Range("W4").Value = 0.8631710361591275
$ws.Range("X4").Value = 0.8639754189990841

# Row 5 updates
$ws.Range("K5").Value = 0.00015
$ws.Range("L5").Value = 0.00005
$ws.Range("Q5").Value = 104
$ws.Range("R5").Value = 115
$ws.Range("T5").Value = 0.1179499228761208
$ws.Range("U5").Value = 0.11790579342659
$ws.Range("W5").Value = 0.863915766140826
$ws.Range("X5").Value = 0.8632329393796765
